# Scheduled Sheets update: refresh cached market-board averages / leve
# profit figures across the eight crafting-profession worksheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 17902.273
$ws.Range("I69").Value = 15404.833
$ws.Range("K69").Value = 46214.499
$ws.Range("M69").Value = -45340.499

$ws.Range("H72").Value = 17902.273
$ws.Range("I72").Value = 15404.833
$ws.Range("K72").Value = 138643.497
$ws.Range("M72").Value = -134275.497

$ws.Range("H82").Value = 2344
$ws.Range("I82").Value = 2344
$ws.Range("K82").Value = 7032
$ws.Range("M82").Value = -6626

$ws.Range("H85").Value = 2344
$ws.Range("I85").Value = 2344
$ws.Range("K85").Value = 7032
$ws.Range("M85").Value = -5628

$ws.Range("H98").Value = 3253.5557
$ws.Range("I98").Value = 2829.2
$ws.Range("J98").Value = 3503.1765
$ws.Range("K98").Value = 2829.2
$ws.Range("L98").Value = 3503.1765
$ws.Range("M98").Value = -1331.2
$ws.Range("N98").Value = -6499.1765

$ws.Range("H113").Value = 2713.6
$ws.Range("I113").Value = 2943.7778
$ws.Range("K113").Value = 2943.7778
$ws.Range("M113").Value = 310.2222000000002

$ws.Range("H122").Value = 3253.5557
$ws.Range("I122").Value = 2829.2
$ws.Range("J122").Value = 3503.1765
$ws.Range("K122").Value = 8487.599999999999
$ws.Range("L122").Value = 10509.5295
$ws.Range("M122").Value = -6037.599999999999
$ws.Range("N122").Value = -15409.5295

$ws.Range("H137").Value = 977.875
$ws.Range("I137").Value = 924.75
$ws.Range("J137").Value = 1243.5
$ws.Range("K137").Value = 2774.25
$ws.Range("L137").Value = 3730.5
$ws.Range("M137").Value = -224.25
$ws.Range("N137").Value = -8830.5


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 59499
$ws.Range("J92").Value = 59499
$ws.Range("L92").Value = 59499
$ws.Range("N92").Value = -64491

$ws.Range("H97").Value = 6319
$ws.Range("I97").Value = 8284.643
$ws.Range("J97").Value = 2879.125
$ws.Range("K97").Value = 8284.643
$ws.Range("L97").Value = 2879.125
$ws.Range("M97").Value = -7788.643
$ws.Range("N97").Value = -3871.125

$ws.Range("H122").Value = 1998.5
$ws.Range("I122").Value = 1998
$ws.Range("K122").Value = 5994
$ws.Range("M122").Value = -3544

$ws.Range("H132").Value = 4273.643
$ws.Range("I132").Value = 4234.1113
$ws.Range("J132").Value = 4344.8
$ws.Range("K132").Value = 12702.3339
$ws.Range("L132").Value = 13034.4
$ws.Range("M132").Value = -10172.3339
$ws.Range("N132").Value = -18094.4


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4639.125
$ws.Range("I94").Value = 2517.1667
$ws.Range("K94").Value = 2517.1667
$ws.Range("M94").Value = -2066.1667

$ws.Range("H99").Value = 6296.846
$ws.Range("I99").Value = 6792.087
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 6792.087
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -5294.087
$ws.Range("N99").Value = -5496

$ws.Range("H134").Value = 2043.4117
$ws.Range("I134").Value = 1844.8148
$ws.Range("J134").Value = 2809.4285
$ws.Range("K134").Value = 5534.4444
$ws.Range("L134").Value = 8428.2855
$ws.Range("M134").Value = -2999.4444
$ws.Range("N134").Value = -13498.2855


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3292.861
$ws.Range("I31").Value = 6377.25
$ws.Range("J31").Value = 2411.6072
$ws.Range("K31").Value = 6377.25
$ws.Range("L31").Value = 2411.6072
$ws.Range("M31").Value = -6082.25
$ws.Range("N31").Value = -3001.6072

$ws.Range("H34").Value = 3292.861
$ws.Range("I34").Value = 6377.25
$ws.Range("J34").Value = 2411.6072
$ws.Range("K34").Value = 6377.25
$ws.Range("L34").Value = 2411.6072
$ws.Range("M34").Value = -6175.25
$ws.Range("N34").Value = -2815.6072


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14898429
$ws.Range("I4").Value = 20474370
$ws.Range("K4").Value = 61423110
$ws.Range("M4").Value = -61422998

$ws.Range("H106").Value = 14247.5
$ws.Range("I106").Value = 7000
$ws.Range("J106").Value = 16663.334
$ws.Range("K106").Value = 21000
$ws.Range("L106").Value = 49990.00199999999
$ws.Range("M106").Value = -20054
$ws.Range("N106").Value = -51882.00199999999

$ws.Range("H132").Value = 2356.5312
$ws.Range("I132").Value = 1918.619
$ws.Range("J132").Value = 3192.5454
$ws.Range("K132").Value = 17267.571
$ws.Range("L132").Value = 28732.9086
$ws.Range("M132").Value = -14737.571
$ws.Range("N132").Value = -33792.9086


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1849.25
$ws.Range("I113").Value = 1499.5
$ws.Range("K113").Value = 1499.5
$ws.Range("M113").Value = 670.5

$ws.Range("H122").Value = 1109.4
$ws.Range("I122").Value = 1151
$ws.Range("J122").Value = 1081.6666
$ws.Range("K122").Value = 3453
$ws.Range("L122").Value = 3244.9998
$ws.Range("M122").Value = -1003
$ws.Range("N122").Value = -8144.9998

$ws.Range("H123").Value = 38114.92
$ws.Range("J123").Value = 37772.184
$ws.Range("L123").Value = 37772.184
$ws.Range("N123").Value = -42672.184

$ws.Range("H126").Value = 4748.6924
$ws.Range("I126").Value = 4498.5
$ws.Range("K126").Value = 13495.5
$ws.Range("M126").Value = -11025.5

$ws.Range("H132").Value = 4903782
$ws.Range("I132").Value = 1955.625
$ws.Range("J132").Value = 9260961
$ws.Range("K132").Value = 5866.875
$ws.Range("L132").Value = 27782883
$ws.Range("M132").Value = -3336.875
$ws.Range("N132").Value = -27787943


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1316.0476
$ws.Range("I16").Value = 1065.8422
$ws.Range("J16").Value = 3693
$ws.Range("K16").Value = 1065.8422
$ws.Range("L16").Value = 3693
$ws.Range("M16").Value = -895.8422
$ws.Range("N16").Value = -4033

$ws.Range("H22").Value = 3057.6667
$ws.Range("I22").Value = 1071.5
$ws.Range("J22").Value = 4381.778
$ws.Range("K22").Value = 1071.5
$ws.Range("L22").Value = 4381.778
$ws.Range("M22").Value = -776.5
$ws.Range("N22").Value = -4971.778

$ws.Range("H27").Value = 3057.6667
$ws.Range("I27").Value = 1071.5
$ws.Range("J27").Value = 4381.778
$ws.Range("K27").Value = 1071.5
$ws.Range("L27").Value = 4381.778
$ws.Range("M27").Value = -964.5
$ws.Range("N27").Value = -4595.778

$ws.Range("H40").Value = 10334.333
$ws.Range("I40").Value = 12751.5
$ws.Range("K40").Value = 12751.5
$ws.Range("M40").Value = -12615.5

$ws.Range("H68").Value = 9999
$ws.Range("J68").Value = 9999
$ws.Range("L68").Value = 9999
$ws.Range("N68").Value = -11497

$ws.Range("H71").Value = 9999
$ws.Range("J71").Value = 9999
$ws.Range("L71").Value = 49995
$ws.Range("N71").Value = -57483

$ws.Range("H82").Value = 2225.5715
$ws.Range("I82").Value = 1668.2858
$ws.Range("J82").Value = 2782.8572
$ws.Range("K82").Value = 1668.2858
$ws.Range("L82").Value = 2782.8572
$ws.Range("M82").Value = -1307.2858
$ws.Range("N82").Value = -3504.8572

$ws.Range("H85").Value = 2225.5715
$ws.Range("I85").Value = 1668.2858
$ws.Range("J85").Value = 2782.8572
$ws.Range("K85").Value = 1668.2858
$ws.Range("L85").Value = 2782.8572
$ws.Range("M85").Value = -420.2858000000001
$ws.Range("N85").Value = -5278.8572

$ws.Range("H122").Value = 2999.5
$ws.Range("I122").Value = 2999.5
$ws.Range("K122").Value = 8998.5
$ws.Range("M122").Value = -6548.5


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 139998
$ws.Range("J76").Value = 139998
$ws.Range("L76").Value = 139998
$ws.Range("N76").Value = -140628

$ws.Range("H79").Value = 139998
$ws.Range("J79").Value = 139998
$ws.Range("L79").Value = 139998
$ws.Range("N79").Value = -142182

$ws.Range("H107").Value = 308.33334
$ws.Range("I107").Value = 297.70834
$ws.Range("J107").Value = 336.66666
$ws.Range("K107").Value = 893.1250200000001
$ws.Range("L107").Value = 1009.99998
$ws.Range("M107").Value = 1026.87498
$ws.Range("N107").Value = -4849.99998

$ws.Range("H122").Value = 1819.3529
$ws.Range("I122").Value = 1819.3529
$ws.Range("K122").Value = 5458.0587
$ws.Range("M122").Value = -3008.0587

$ws.Range("H132").Value = 1715.9706
$ws.Range("I132").Value = 1572.9656
$ws.Range("J132").Value = 2545.4
$ws.Range("K132").Value = 4718.8968
$ws.Range("L132").Value = 7636.200000000001
$ws.Range("M132").Value = -2188.8968
$ws.Range("N132").Value = -12696.2

$ws.Range("H135").Value = 82345.5
$ws.Range("J135").Value = 82345.5
$ws.Range("L135").Value = 82345.5
$ws.Range("N135").Value = -92485.5

$ws.Range("H139").Value = 199990
$ws.Range("J139").Value = 199990
$ws.Range("L139").Value = 199990
$ws.Range("N139").Value = -210270

